# Adding more words from today's lessons - up to 252!
# Appends 10 new Farsi/Pronunciation/English vocabulary rows (244-253)
# to the end of the word list on Sheet1 (columns A:C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Farsi word, Pronunciation (optional, may be blank), English meaning
$rows = @(
    @("انجام", "anjam", "do"),
    @("قرار", "gherar", "appointment"),
    @("محل", "Mahl", "location"),
    @("فوق ", "", "Above"),
    @("العاده", "", "Extraordinary"),
    @("می گردیم", "", "We return"),
    @("تشنه", "Tashneh", "Thirsty"),
    @("واقع", "", "Indeed"),
    @("بخوانم", "", "I read"),
    @("نام", "Nam(eh)", "name")
)

$startRow = 244
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $entry = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $entry[0]
    if ($entry[1] -ne "") {
        $ws.Cells.Item($r, 2).Value = $entry[1]
    }
    $ws.Cells.Item($r, 3).Value = $entry[2]
}

# Match the author's final cursor position / selection from the edit session.
$ws.Range("B252").Select() | Out-Null
